$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F column values for rows 1-25 that flipped between blank/number.
# (column F holds header "F" in row 1; rows 2..25 are data rows "RM ...")
$ws.Range("F6").Value = 16.43
$ws.Range("F8").Value = ""
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = ""
$ws.Range("F17").Value = 17.78
$ws.Range("F18").Value = 18.35
$ws.Range("F19").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("F23").Value = 16.48

# Remove the two rows that were dropped from the bottom block ("RM 232" and "SC 92").
# Delete the lower one first so row indices above it stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the deletions, rows 26-33 hold (in order): SC 5, SC 101, SC 105, SC 119,
# SC 120, SC 132, SC 193, SC 232. Apply the remaining value changes to them.
$ws.Range("D27").Value = -14.6
$ws.Range("F27").Value = ""

$ws.Range("D29").Value = ""

$ws.Range("D32").Value = ""
